$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.057720695672616
$arr[0,2] = 1.058854182357834
$arr[0,3] = 1.070054411905201
$arr[0,4] = 1.074568326332472
$ws.Range("B2:F2").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040459682596729
$arr[0,1] = 1.062715647464937
$arr[0,2] = 1.061585279934077
$arr[0,3] = 1.072755273275143
$arr[0,4] = 1.077257194542978
$arr[0,5] = 1.024530879045038
$ws.Range("I2:N2").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.059067513419303
$arr[0,2] = 1.059878753034176
$arr[0,3] = 1.071345258861037
$arr[0,4] = 1.075835263038431
$ws.Range("B3:F3").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040715595236939
$arr[0,1] = 1.063713136872082
$arr[0,2] = 1.062423432336208
$arr[0,3] = 1.073861233716562
$arr[0,4] = 1.078340177303941
$arr[0,5] = 1.024873385616457
$ws.Range("I3:N3").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.059938565262109
$arr[0,2] = 1.060541114971
$arr[0,3] = 1.07218041658044
$arr[0,4] = 1.076654888986313
$ws.Range("B4:F4").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040879559990413
$arr[0,1] = 1.064357665112678
$arr[0,2] = 1.062964561864413
$arr[0,3] = 1.074576201175091
$arr[0,4] = 1.079040210142609
$arr[0,5] = 1.025094445407024
$ws.Range("I4:N4").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.060304656254151
$arr[0,2] = 1.060819429213237
$arr[0,3] = 1.072531494113114
$arr[0,4] = 1.076999422380611
$ws.Range("B5:F5").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040948101666878
$arr[0,1] = 1.064628408263548
$arr[0,2] = 1.06319176456648
$arr[0,3] = 1.074876617255552
$arr[0,4] = 1.079334331655579
$arr[0,5] = 1.025187244379746
$ws.Range("I5:N5").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.060366118828998
$arr[0,2] = 1.060866151083412
$arr[0,3] = 1.072590440319843
$arr[0,4] = 1.07705726893799
$ws.Range("B6:F6").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040959587316027
$arr[0,1] = 1.064673854600206
$arr[0,2] = 1.063229896054542
$arr[0,3] = 1.074927049346804
$arr[0,4] = 1.079383705916177
$arr[0,5] = 1.025202817864053
$ws.Range("I6:N6").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.059943457373208
$arr[0,2] = 1.060544834378026
$arr[0,3] = 1.072185107781678
$arr[0,4] = 1.076659492801951
$ws.Range("B7:F7").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040880477375824
$arr[0,1] = 1.064361283643812
$arr[0,2] = 1.062967598888053
$arr[0,3] = 1.074580215960313
$arr[0,4] = 1.079044140882751
$arr[0,5] = 1.025095685919352
$ws.Range("I7:N7").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.058175950077262
$arr[0,2] = 1.059200566943459
$arr[0,3] = 1.070490683398212
$arr[0,4] = 1.074996529621768
$ws.Range("B8:F8").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.040546506889252
$arr[0,1] = 1.063052944349056
$arr[0,2] = 1.061868789315034
$arr[0,3] = 1.073129176525717
$arr[0,4] = 1.077623346076555
$arr[0,5] = 1.024646748105456
$ws.Range("I8:N8").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.05505792825001
$arr[0,2] = 1.056827088769576
$arr[0,3] = 1.067503951635104
$arr[0,4] = 1.072064782258586
$ws.Range("B9:F9").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039945513696431
$arr[0,1] = 1.060740378187661
$arr[0,2] = 1.059923189901885
$arr[0,3] = 1.070567071983935
$arr[0,4] = 1.07511403260705
$arr[0,5] = 1.023851306035132
$ws.Range("I9:N9").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.052976702765791
$arr[0,2] = 1.05524148679772
$arr[0,3] = 1.065511990765251
$arr[0,4] = 1.070109182082043
$ws.Range("B10:F10").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039536411664961
$arr[0,1] = 1.059193742530816
$arr[0,2] = 1.058619717664229
$arr[0,3] = 1.068855357484336
$arr[0,4] = 1.073437181260755
$arr[0,5] = 1.023318039298583
$ws.Range("I10:N10").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.052074854823109
$arr[0,2] = 1.054554097782034
$arr[0,3] = 1.064649218856421
$arr[0,4] = 1.069262086807643
$ws.Range("B11:F11").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039357255002106
$arr[0,1] = 1.058522830566315
$arr[0,2] = 1.058053754086553
$arr[0,3] = 1.068113266100116
$arr[0,4] = 1.07271010924588
$arr[0,5] = 1.023086414079067
$ws.Range("I11:N11").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.051739764180498
$arr[0,2] = 1.054298646330393
$arr[0,3] = 1.06432870743692
$arr[0,4] = 1.068947388386919
$ws.Range("B12:F12").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039290404993391
$arr[0,1] = 1.05827343944499
$arr[0,2] = 1.057843294560423
$arr[0,3] = 1.067837480784674
$arr[0,4] = 1.072439891429602
$arr[0,5] = 1.023000269511232
$ws.Range("I12:N12").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.051811647046302
$arr[0,2] = 1.054353447210037
$arr[0,3] = 1.064397460069108
$arr[0,4] = 1.069014894560776
$ws.Range("B13:F13").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.03930475828345
$arr[0,1] = 1.058326943086145
$arr[0,2] = 1.057888449535995
$arr[0,3] = 1.067896644036116
$arr[0,4] = 1.072497860945445
$arr[0,5] = 1.023018752744474
$ws.Range("I13:N13").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.052047158276835
$arr[0,2] = 1.054532984640592
$arr[0,3] = 1.064622726108979
$arr[0,4] = 1.069236074767457
$ws.Range("B14:F14").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039351735350612
$arr[0,1] = 1.058502219600381
$arr[0,2] = 1.058036362253349
$arr[0,3] = 1.068090472464692
$arr[0,4] = 1.072687776052637
$arr[0,5] = 1.023079295563813
$ws.Range("I14:N14").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.052192250582411
$arr[0,2] = 1.054643586997806
$arr[0,3] = 1.064761514609279
$arr[0,4] = 1.069372344531288
$ws.Range("B15:F15").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039380639263449
$arr[0,1] = 1.058610188702249
$arr[0,2] = 1.058127464889542
$arr[0,3] = 1.068209877972308
$arr[0,4] = 1.072804768932407
$arr[0,5] = 1.023116583569118
$ws.Range("I15:N15").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053036541081926
$arr[0,2] = 1.055287089208365
$arr[0,3] = 1.065569244746274
$arr[0,4] = 1.070165394243898
$ws.Range("B16:F16").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039548259223226
$arr[0,1] = 1.059238243032244
$arr[0,2] = 1.058657245913601
$arr[0,3] = 1.068904588251651
$arr[0,4] = 1.073485413640432
$arr[0,5] = 1.023333396314521
$ws.Range("I16:N16").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053565961399872
$arr[0,2] = 1.055690521791932
$arr[0,3] = 1.066075845875404
$arr[0,4] = 1.070662768711466
$ws.Range("B17:F17").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039652863318022
$arr[0,1] = 1.05963187926644
$arr[0,2] = 1.05898914651168
$arr[0,3] = 1.069340116408976
$arr[0,4] = 1.073912098429754
$arr[0,5] = 1.023469204574564
$ws.Range("I17:N17").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053874699592638
$arr[0,2] = 1.055925758925853
$arr[0,3] = 1.066371314853392
$arr[0,4] = 1.070952849381192
$ws.Range("B18:F18").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039713682970384
$arr[0,1] = 1.05986136422038
$arr[0,2] = 1.059182588979131
$arr[0,3] = 1.069594065248964
$arr[0,4] = 1.074160881681578
$arr[0,5] = 1.023548350096455
$ws.Range("I18:N18").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053979960664031
$arr[0,2] = 1.056005955507151
$arr[0,3] = 1.066472058362122
$arr[0,4] = 1.071051754499071
$ws.Range("B19:F19").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039734387994966
$arr[0,1] = 1.05993959300504
$arr[0,2] = 1.059248522568426
$arr[0,3] = 1.069680640517025
$arr[0,4] = 1.074245694357776
$arr[0,5] = 1.023575324973905
$ws.Range("I19:N19").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053509166234939
$arr[0,2] = 1.055647245407924
$arr[0,3] = 1.066021494762915
$arr[0,4] = 1.070609408208867
$ws.Range("B20:F20").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.03964166037664
$arr[0,1] = 1.059589657886479
$arr[0,2] = 1.05895355220383
$arr[0,3] = 1.069293397420181
$arr[0,4] = 1.073866329015119
$arr[0,5] = 1.023454640791731
$ws.Range("I20:N20").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.051977809000751
$arr[0,2] = 1.054480118800395
$arr[0,3] = 1.064556391993951
$arr[0,4] = 1.069170944103081
$ws.Range("B21:F21").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039337910153156
$arr[0,1] = 1.05845061015394
$arr[0,2] = 1.057992812163402
$arr[0,3] = 1.068033398700792
$arr[0,4] = 1.07263185497686
$arr[0,5] = 1.023061470218957
$ws.Range("I21:N21").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.051014379054532
$arr[0,2] = 1.053745578066215
$arr[0,3] = 1.06363499334932
$arr[0,4] = 1.068266236469808
$ws.Range("B22:F22").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.03914517541317
$arr[0,1] = 1.057733376428005
$arr[0,2] = 1.057387393327719
$arr[0,3] = 1.067240378621795
$arr[0,4] = 1.071854817605509
$arr[0,5] = 1.022813639221976
$ws.Range("I22:N22").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.051525170328086
$arr[0,2] = 1.054135041252944
$arr[0,3] = 1.064123467215509
$arr[0,4] = 1.068745867530245
$ws.Range("B23:F23").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039247514376036
$arr[0,1] = 1.058113697942804
$arr[0,2] = 1.057708467336119
$arr[0,3] = 1.067660851300496
$arr[0,4] = 1.07226682372473
$arr[0,5] = 1.022945079050099
$ws.Range("I23:N23").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.053534829721946
$arr[0,2] = 1.05566680038485
$arr[0,3] = 1.066046053761759
$arr[0,4] = 1.070633519614006
$ws.Range("B24:F24").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.039646723103499
$arr[0,1] = 1.059608736270251
$arr[0,2] = 1.058969636201394
$arr[0,3] = 1.069314507990073
$arr[0,4] = 1.073887010536203
$arr[0,5] = 1.023461221751918
$ws.Range("I24:N24").Value = $arr

$arr = New-Object 'object[,]' 1,5
$arr[0,0] = 1.02
$arr[0,1] = 1.055864443465935
$arr[0,2] = 1.057441259734685
$arr[0,3] = 1.068276222661767
$arr[0,4] = 1.072822891318446
$ws.Range("B25:F25").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 1.04010236953703
$arr[0,1] = 1.061339089142802
$arr[0,2] = 1.060427294601189
$arr[0,3] = 1.07123006829704
$arr[0,4] = 1.075763440060588
$arr[0,5] = 1.024057467538337
$ws.Range("I25:N25").Value = $arr
